# Add the quarterly "2022-Q1" fund-holdings sheet and fold its summary row
# into the "总计" (Total) sheet.
#
# Layout before the edit:  2021-Q2 | 2021-Q3 | 2021-Q4 | 总计
# Layout after the edit:   2021-Q2 | 2021-Q3 | 2021-Q4 | 2022-Q1 | 总计
#
# The existing "总计" sheet keeps its physical identity/sheetId and is
# renamed to "2022-Q1", then repopulated with the fund-holdings table for
# the new quarter. A fresh duplicate becomes the new "总计" sheet, which
# gets the old total-table content (already there from the duplication)
# plus one new row summarising 2022-Q1.

function Set-TextValue {
    # Writes $value into $range while forcing Excel to keep it as literal
    # text (so numeric-looking strings like "82.26" or "014133" are not
    # silently reinterpreted as numbers), without leaving a stray
    # "number stored as text" style behind.
    param($range, $value)
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

$wb = $excel.ActiveWorkbook

# --- Duplicate "总计" so we have a correctly formatted sheet to turn into
#     the new "总计", while the in-place original becomes "2022-Q1" (this
#     is what keeps the sheetId/r:id numbering lined up with the diff).
$totalOrig = $wb.Worksheets.Item("总计")
$totalOrig.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$totalOrig.Name = "2022-Q1"
$total = $wb.Worksheets.Item("总计 (2)")
$total.Name = "总计"

$q1 = $wb.Worksheets.Item("2022-Q1")
$q4 = $wb.Worksheets.Item("2021-Q4")

# --- Bring over the E:H header formatting from the 2021-Q4 fund-holdings
#     sheet (same 8-column layout we need here).
$q4.Range("E1:H1").Copy($q1.Range("E1"))

# Header row text (reusing formatting already present on B1:D1 / E1:H1).
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# --- Fund holdings data (rows 2-4). Columns B-G are text (even the
#     numeric-looking ones); H is a genuine number.
Set-TextValue $q1.Range("B2") "161017"
$q1.Range("C2").Value = "富国中证500指数增强(LOF)"
Set-TextValue $q1.Range("D2") "82.26"
Set-TextValue $q1.Range("E2") "91.34"
Set-TextValue $q1.Range("F2") "0.96"
Set-TextValue $q1.Range("G2") "0.7897"
$q1.Range("H2").Value = 5

Set-TextValue $q1.Range("B3") "014133"
$q1.Range("C3").Value = "工银中证500六个月持有指数增强A"
Set-TextValue $q1.Range("D3") "3.07"
Set-TextValue $q1.Range("E3") "93.69"
Set-TextValue $q1.Range("F3") "0.81"
Set-TextValue $q1.Range("G3") "0.0249"
$q1.Range("H3").Value = 8

Set-TextValue $q1.Range("B4") "014134"
$q1.Range("C4").Value = "工银中证500六个月持有指数增强C"
Set-TextValue $q1.Range("D4") "1.12"
Set-TextValue $q1.Range("E4") "93.69"
Set-TextValue $q1.Range("F4") "0.81"
Set-TextValue $q1.Range("G4") "0.0091"
$q1.Range("H4").Value = 8

# --- Insert the new 2022-Q1 summary row into "总计" and renumber the
#     index column (A) that follows it.
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()
$total.Range("A3").Copy($total.Range("A2"))

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.82

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# Restore the original active sheet/selection state (the "总计" duplication
# leaves itself active/selected, which the source workbook wasn't).
$wb.Worksheets.Item(1).Activate()
